# Applies the "Updated cryptos list" data refresh (prices + 1h volume %)
# to Sheet1 of the cryptos workbook. Two rows (33/34, 40/41, 46/47) had their
# coin/link/price/volume swapped with their neighbour as the ranking shifted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.788.42'
$ws.Range("E2").Value = '  -0.59%  '

# Row 3
$ws.Range("D3").Value = '1.906.03'
$ws.Range("E3").Value = '  +0.12%  '

# Row 4
$ws.Range("E4").Value = '  -0.41%  '

# Row 5
$ws.Range("D5").Value = '''313.05'
$ws.Range("E5").Value = '  -1.25%  '

# Row 6
$ws.Range("D6").Value = '''1.002'

# Row 7
$ws.Range("D7").Value = '''0.4978'
$ws.Range("E7").Value = '  +3.28%  '

# Row 8
$ws.Range("D8").Value = '''0.3775'
$ws.Range("E8").Value = '  -0.50%  '

# Row 9
$ws.Range("D9").Value = '''0.07253'
$ws.Range("E9").Value = '  -1.55%  '

# Row 10
$ws.Range("D10").Value = '''21.11'
$ws.Range("E10").Value = '  +1.69%  '

# Row 11
$ws.Range("D11").Value = '''0.9016'
$ws.Range("E11").Value = '  -3.19%  '

# Row 12
$ws.Range("D12").Value = '''0.07630'
$ws.Range("E12").Value = '  -1.45%  '

# Row 13
$ws.Range("D13").Value = '1.870.30'
$ws.Range("E13").Value = '  -1.77%  '

# Row 14
$ws.Range("D14").Value = '''5.457'
$ws.Range("E14").Value = '  -0.44%  '

# Row 15
$ws.Range("D15").Value = '''91.81'
$ws.Range("E15").Value = '  +0.14%  '

# Row 16
$ws.Range("E16").Value = '  -0.40%  '

# Row 17
$ws.Range("D17").Value = '''0.000008699'
$ws.Range("E17").Value = '  -1.92%  '

# Row 18
$ws.Range("E18").Value = '  -0.39%  '

# Row 19
$ws.Range("D19").Value = '27.831.92'
$ws.Range("E19").Value = '  -0.56%  '

# Row 20
$ws.Range("D20").Value = '''14.53'
$ws.Range("E20").Value = '  -0.88%  '

# Row 21
$ws.Range("D21").Value = '''5.159'
$ws.Range("E21").Value = '  +0.23%  '

# Row 22
$ws.Range("D22").Value = '2.132.85'
$ws.Range("E22").Value = '  -1.24%  '

# Row 23
$ws.Range("D23").Value = '''10.81'
$ws.Range("E23").Value = '  -0.77%  '

# Row 24
$ws.Range("D24").Value = '''6.583'
$ws.Range("E24").Value = '  -0.62%  '

# Row 25
$ws.Range("D25").Value = '''153.07'
$ws.Range("E25").Value = '  -1.91%  '

# Row 26
$ws.Range("D26").Value = '''1.846'
$ws.Range("E26").Value = '  -3.94%  '

# Row 27
$ws.Range("D27").Value = '''2.210'
$ws.Range("E27").Value = '  +3.91%  '

# Row 28
$ws.Range("D28").Value = '''18.34'
$ws.Range("E28").Value = '  -0.70%  '

# Row 29
$ws.Range("D29").Value = '''114.99'
$ws.Range("E29").Value = '  -1.76%  '

# Row 30
$ws.Range("D30").Value = '''4.864'
$ws.Range("E30").Value = '  -1.97%  '

# Row 31
$ws.Range("D31").Value = '''0.08944'
$ws.Range("E31").Value = '  +0.10%  '

# Row 32
$ws.Range("D32").Value = '''3.193'
$ws.Range("E32").Value = '  -2.08%  '

# Row 33
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7865'
$ws.Range("E33").Value = '  +2.58%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''1.232'
$ws.Range("E34").Value = '  -1.63%  '

# Row 35
$ws.Range("D35").Value = '''4.785'
$ws.Range("E35").Value = '  +2.54%  '

# Row 36
$ws.Range("D36").Value = '''2.635'
$ws.Range("E36").Value = '  +3.37%  '

# Row 37
$ws.Range("D37").Value = '''0.02075'
$ws.Range("E37").Value = '  +0.91%  '

# Row 38
$ws.Range("D38").Value = '''3.057'
$ws.Range("E38").Value = '  +2.05%  '

# Row 39
$ws.Range("E39").Value = '  -1.15%  '

# Row 40
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '''0.05297'
$ws.Range("E40").Value = '  +0.42%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.5502'
$ws.Range("E41").Value = '  +0.30%  '

# Row 42
$ws.Range("D42").Value = '''6.755'
$ws.Range("E42").Value = '  -2.63%  '

# Row 43
$ws.Range("D43").Value = '''113.99'
$ws.Range("E43").Value = '  +3.74%  '

# Row 44
$ws.Range("D44").Value = '''8.451'
$ws.Range("E44").Value = '  -0.31%  '

# Row 45
$ws.Range("D45").Value = '''0.1507'
$ws.Range("E45").Value = '  -1.21%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''10.52'
$ws.Range("E46").Value = '  -1.52%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.4780'
$ws.Range("E47").Value = '  -0.51%  '

# Row 48
$ws.Range("E48").Value = '  -0.37%  '

# Row 49
$ws.Range("D49").Value = '''1.628'
$ws.Range("E49").Value = '  -0.98%  '

# Row 50
$ws.Range("D50").Value = '''67.16'
$ws.Range("E50").Value = '  -0.96%  '

# Row 51
$ws.Range("D51").Value = '''0.06023'
$ws.Range("E51").Value = '  -0.91%  '
